$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data for row 15 (previously blank placeholder cells)
$ws.Range("A15:D15").HorizontalAlignment = -4108
$ws.Range("A15").Value = "Hj Nia Herdiani, SE., M.pd"
$ws.Range("B15").Value = "MDK"
$ws.Range("C15").Value = 907
$ws.Range("D15").Value = "Test"

# Update the active selection to match the saved view state
$ws.Range("E11").Select()
